# Moved to version 0.14.17
# Adds a new day (row 13) of coverage data to Sheet1, which grows the
# line-chart source ranges from rows 2:12 to 2:13, and repositions the
# chart on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the new data row (row 13) -----------------------------------
# Carry the formatting down from row 12 (date / 2-decimal number formats).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N12:R12").Copy()
$ws.Range("N13:R13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 44084
$ws.Range("B13").Value = 12

$ws.Range("D13").Value = 161
$ws.Range("E13").Value = 87
$ws.Range("F13").Value = 1456
$ws.Range("G13").Value = 580
$ws.Range("H13").Value = 47
$ws.Range("I13").Value = 20
$ws.Range("J13").Value = 14
$ws.Range("K13").Value = 210
$ws.Range("L13").Value = 83

$ws.Range("N13").Formula = "=100*E13/D13"
$ws.Range("O13").Formula = "=100*G13/F13"
$ws.Range("P13").Formula = "=100*H13/D13"
$ws.Range("Q13").Formula = "=100*J13/I13"
$ws.Range("R13").Formula = "=100*L13/K13"

# --- Extend the chart series ranges to include the new row -----------------
$chart = $ws.ChartObjects(1).Chart
$valCols = @("N", "O", "P", "Q", "R")
for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
    $col = $valCols[$i - 1]
    $series = $chart.SeriesCollection($i)
    $series.Formula = "=SERIES(Sheet1!`$$col`$1,Sheet1!`$B`$2:`$B`$13,Sheet1!`$$col`$2:`$$col`$13,$i)"
}

# --- Reposition the chart ---------------------------------------------------
# The chart keeps its size but is nudged down/right on the sheet (to clear
# the newly added data row). Translate by the same amount in points
# (30pt right, 39.75pt down).
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = $chartObj.Left + 30
$chartObj.Top = $chartObj.Top + 39.75

# --- Update the selected cell ----------------------------------------------
$ws.Range("P29").Select()
